$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9975482225418091
$ws.Range("B1").Value = 2.239984273910522
$ws.Range("C1").Value = 4.935811519622803
$ws.Range("D1").Value = 1.695178866386414
$ws.Range("E1").Value = 1.289914608001709
